# Auto-generated Excel COM-interop script applying the Aegis_Profits.xlsx diff.
# For each affected row, update the changed numeric cells (H-N), adding or
# clearing cells where the diff adds/removes them entirely.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 725.1818
$ws.Range("I32").Value = 716.6667
$ws.Range("J32").Value = 728.375
$ws.Range("K32").Value = 716.6667
$ws.Range("L32").Value = 728.375
$ws.Range("M32").Value = -390.6667
$ws.Range("N32").Value = -1380.375

$ws.Range("H51").Value = 24159.8
$ws.Range("J51").Value = 6749.5
$ws.Range("L51").Value = 6749.5
$ws.Range("N51").Value = -7717.5

$ws.Range("H116").Value = 4995
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 4993.3335
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 4993.3335
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -11877.3335

$ws.Range("H128").Value = 45446.332
$ws.Range("J128").Value = 45446.332
$ws.Range("L128").Value = 45446.332
$ws.Range("N128").Value = -55406.332

$ws.Range("H129").Value = 868.1896400000001
$ws.Range("J129").Value = 969
$ws.Range("L129").Value = 2907
$ws.Range("N129").Value = -12907

$ws.Range("H138").Value = 2705.328
$ws.Range("I138").Value = 1595.1904
$ws.Range("J138").Value = 3288.15
$ws.Range("K138").Value = 4785.5712
$ws.Range("L138").Value = 9864.450000000001
$ws.Range("M138").Value = 354.4287999999997
$ws.Range("N138").Value = -20144.45

$ws.Range("H141").Value = 4210.3335
$ws.Range("I141").Value = 2429.875
$ws.Range("J141").Value = 7771.25
$ws.Range("K141").Value = 7289.625
$ws.Range("L141").Value = 23313.75
$ws.Range("M141").Value = -2109.625
$ws.Range("N141").Value = -33673.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 9000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 9000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = -9346

$ws.Range("H9").Value = 13999.667
$ws.Range("I9").Value = 13999
$ws.Range("K9").Value = 13999
$ws.Range("M9").Value = -13829

$ws.Range("H20").Value = 13999.667
$ws.Range("I20").Value = 13999
$ws.Range("K20").Value = 13999
$ws.Range("M20").Value = -13729

$ws.Range("H23").Value = 55004.332
$ws.Range("J23").Value = 47503.5
$ws.Range("L23").Value = 47503.5
$ws.Range("N23").Value = -48021.5

$ws.Range("H41").Value = 4085.3333
$ws.Range("I41").Value = 3302.4
$ws.Range("J41").Value = 8000
$ws.Range("K41").Value = 3302.4
$ws.Range("L41").Value = 8000
$ws.Range("M41").Value = -2888.4
$ws.Range("N41").Value = -8828

$ws.Range("H55").Value = 14150
$ws.Range("J55").Value = 14885.714
$ws.Range("L55").Value = 14885.714
$ws.Range("N55").Value = -15515.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1738
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1738
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1738
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -2640

$ws.Range("H122").Value = 2743.7856
$ws.Range("I122").Value = 2908.3333
$ws.Range("J122").Value = 1756.5
$ws.Range("K122").Value = 8724.999899999999
$ws.Range("L122").Value = 5269.5
$ws.Range("M122").Value = -6274.999899999999
$ws.Range("N122").Value = -10169.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3300
$ws.Range("J39").Value = 3300
$ws.Range("L39").Value = 9900
$ws.Range("N39").Value = -10488

$ws.Range("H55").Value = 4612.5
$ws.Range("J55").Value = 4612.5
$ws.Range("L55").Value = 13837.5
$ws.Range("N55").Value = -14191.5

$ws.Range("H68").Value = 16137.449
$ws.Range("J68").Value = 20310.277
$ws.Range("L68").Value = 60930.83099999999
$ws.Range("N68").Value = -62552.83099999999

$ws.Range("H71").Value = 16137.449
$ws.Range("J71").Value = 20310.277
$ws.Range("L71").Value = 182792.493
$ws.Range("N71").Value = -190904.493

$ws.Range("H131").Value = 879.8953
$ws.Range("I131").Value = 440
$ws.Range("J131").Value = 925.0128
$ws.Range("K131").Value = 1320
$ws.Range("L131").Value = 2775.0384
$ws.Range("M131").Value = 3720
$ws.Range("N131").Value = -12855.0384

$ws.Range("H133").Value = 4372.222
$ws.Range("I133").Value = 729.1667
$ws.Range("J133").Value = 7286.6665
$ws.Range("K133").Value = 2187.5001
$ws.Range("L133").Value = 21859.9995
$ws.Range("M133").Value = 2872.4999
$ws.Range("N133").Value = -31979.9995

$ws.Range("H136").Value = 3676.4
$ws.Range("I136").Value = 363.84616
$ws.Range("J136").Value = 7265
$ws.Range("K136").Value = 1091.53848
$ws.Range("L136").Value = 21795
$ws.Range("M136").Value = 4008.46152
$ws.Range("N136").Value = -31995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5698.75
$ws.Range("I43").Value = 1863.3334
$ws.Range("J43").Value = 8000
$ws.Range("K43").Value = 1863.3334
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = -1712.3334
$ws.Range("N43").Value = -8302

$ws.Range("H49").Value = 12240
$ws.Range("J49").Value = 12240
$ws.Range("L49").Value = 12240
$ws.Range("N49").Value = -12608

$ws.Range("H55").Value = 5030
$ws.Range("I55").Value = 5030
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 5030
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -4703
$ws.Range("N55").Value = ""

$ws.Range("H70").Value = 185973.19
$ws.Range("I70").Value = 253963.38
$ws.Range("J70").Value = 4666
$ws.Range("K70").Value = 253963.38
$ws.Range("L70").Value = 4666
$ws.Range("M70").Value = -253693.38
$ws.Range("N70").Value = -5206

$ws.Range("H73").Value = 185973.19
$ws.Range("I73").Value = 253963.38
$ws.Range("J73").Value = 4666
$ws.Range("K73").Value = 253963.38
$ws.Range("L73").Value = 4666
$ws.Range("M73").Value = -253027.38
$ws.Range("N73").Value = -6538

$ws.Range("H132").Value = 3217.6428
$ws.Range("I132").Value = 2692.7778
$ws.Range("J132").Value = 4162.4
$ws.Range("K132").Value = 8078.3334
$ws.Range("L132").Value = 12487.2
$ws.Range("M132").Value = -5548.3334
$ws.Range("N132").Value = -17547.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 391.72415
$ws.Range("I22").Value = 334
$ws.Range("J22").Value = 422.10526
$ws.Range("K22").Value = 334
$ws.Range("L22").Value = 422.10526
$ws.Range("M22").Value = -39
$ws.Range("N22").Value = -1012.10526

$ws.Range("H27").Value = 391.72415
$ws.Range("I27").Value = 334
$ws.Range("J27").Value = 422.10526
$ws.Range("K27").Value = 334
$ws.Range("L27").Value = 422.10526
$ws.Range("M27").Value = -227
$ws.Range("N27").Value = -636.10526

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""

$ws.Range("H122").Value = 2722.5557
$ws.Range("I122").Value = 3167.3333
$ws.Range("J122").Value = 1833
$ws.Range("K122").Value = 9501.999899999999
$ws.Range("L122").Value = 5499
$ws.Range("M122").Value = -7051.999899999999
$ws.Range("N122").Value = -10399
